$wb = $excel.ActiveWorkbook

# Sheets (1-based): 1 = 2EXT04_DNA, 2 = SwateTemplateMetadata, 3 = SRA_GENOMICS
$wsMeta = $wb.Worksheets.Item(2)

# Rename the metadata sheet to its new name.
$wsMeta.Name = "isa_template"

# Remove the ER tags (SRA / GENBANK) from this non-ER template.
$wsMeta.Range("B8").Value = $null
$wsMeta.Range("C8").Value = $null

# Replace the non-ER "Extraction" tag with the proper lower-case "extraction"
# term and give it a Term Accession Number (the Term Source REF column, C12,
# already correctly points at "DNA" and is left untouched).
$wsMeta.Range("B12").Value = "extraction"
$wsMeta.Range("B13").Value = "http://purl.obolibrary.org/obo/OBI_0302884"

# Move the active tab / selection from SRA_GENOMICS to isa_template.
$wsMeta.Activate()
$wsMeta.Range("B15").Select() | Out-Null
